# Applies the "todos los ataques mas docs" edit:
#   - splits " rho factorizationg" into " rho " + spell-checked "factorizationg"
#   - appends a blank paragraph, a "Paradoja del cumpleaños" paragraph,
#     another blank paragraph, and a final paragraph with a trailing code
#     comment / Java expression, moving the _GoBack bookmark into that
#     last paragraph.

$d = $word.ActiveDocument
$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# Locate the paragraph that currently ends in "... Pollar rho factorizationg"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "*rho factorizationg*") {
        $target = $d.Paragraphs.Item($i)
        break
    }
}

# --- Step 1: split the trailing word off into its own spell-checked run ---
$full = $target.Range.Text
$idx = $full.IndexOf("factorizationg")
$wordStart = $target.Range.Start + $idx
$wordEnd = $wordStart + [string]"factorizationg".Length
$wordRange = $d.Range($wordStart, $wordEnd)
$splitXml = "<w:p $wNs>" + `
              "<w:proofErr w:type='spellStart'/>" + `
              "<w:r><w:t>factorizationg</w:t></w:r>" + `
              "<w:proofErr w:type='spellEnd'/>" + `
            "</w:p>"
$null = $wordRange.InsertXML($splitXml)

# --- Step 2: the _GoBack bookmark is about to move to the new last
#     paragraph, so drop the old one before re-adding it further down ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $null = $d.Bookmarks.Item("_GoBack").Delete()
}

# --- Step 3: append the four new paragraphs after the (now bookmark-free)
#     target paragraph ---
$insertPos = $target.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)

$cumpleanos = "Paradoja del cumplea" + [char]0x00F1 + "os"

$newParasXml = "<w:p $wNs/>" + `
               "<w:p $wNs><w:r><w:t>$cumpleanos</w:t></w:r></w:p>" + `
               "<w:p $wNs/>" + `
               "<w:p $wNs>" + `
                 "<w:r><w:t>// comprobar, pero no sirve para</w:t></w:r>" + `
                 "<w:r><w:t xml:space='preserve'> nada</w:t></w:r>" + `
                 "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" + `
                 "<w:bookmarkEnd w:id='0'/>" + `
                 "<w:r><w:t xml:space='preserve'> ((w.multiply(s)).add(publica.multiply(t))).equals(BigInteger.ONE));</w:t></w:r>" + `
               "</w:p>"

$null = $insertRange.InsertXML($newParasXml)
